$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.26%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.75%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.685"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.59%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08350"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.45%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.808"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.70%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.17%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.48%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.896"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.81%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9263"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.36%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1286"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.19%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1981"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.93%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09499"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.80%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03847"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1061"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.96%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001305"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.24%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.05%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.08%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.55%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.664"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.01%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.89%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04418"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.30%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001277"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.14%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004389"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.85%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.93%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02820"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.50%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05534"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.82%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007946"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.03%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1432"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.97%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009301"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.29%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002111"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.62%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.99%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006932"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.58%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "14.12%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.14%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"
